$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 data mirrors row 2's structure, for a new AT-003 site entry.
$ws.Range("A3").Value = "icmp"
$ws.Range("B3").Value = "icmp-at003_sdwan_at-003_130.143.164.190"
$ws.Range("C3").Value = "icmp-at003_sdwan_at-003_130.143.164.190"
$ws.Range("D3").Value = "130.143.164.200"
$ws.Range("E3").Value = $true
$ws.Range("F3").Value = $false
$ws.Range("G3").Value = "Klagenfurt"
$ws.Range("H3").Value = "AT"
$ws.Range("I3").Value = "Austria"
$ws.Range("J3").Value = 46.61
$ws.Range("K3").Value = 14.32
$ws.Range("L3").Value = "Klagenfurt Koningsbergerstrasse"
$ws.Range("M3").Value = "abd1aff66f13420040ae0d55eb3ee46f"
$ws.Range("N3").Value = "AT-003"
$ws.Range("O3").Value = "AT-003_Philips PH"
$ws.Range("P3").Value = "1bcee5e1dbb6f3844bdbee71ca961979"
$ws.Range("R3").Value = "at003_sdwan_at-003_130.143.164.190"
$ws.Range("T3").Value = "AT-003-SITE-UPLINK"
$ws.Range("V3").Value = "Application.Bridge.Monitoring"
$ws.Range("W3").Value = "any"
$ws.Range("X3").Value = "16s"
$ws.Range("Y3").Value = "1s"
$ws.Range("Z3").Value = "site"

# Update selection to B6 (also resets the view's top-left scroll position)
$ws.Range("B6").Select() | Out-Null
